# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets. Rows 3,4,5,6,7,8,10,13,15,19,20,21,22 get new counts.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 306
    4  = 877
    5  = 43
    6  = 343
    7  = 10689
    8  = 274
    10 = 6
    13 = 138
    15 = 45
    19 = 301
    20 = 1052
    21 = 50
    22 = 104
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
